$p = $ppt.ActivePresentation
$newDate = "27-06-2023"

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tf = $shape.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            if ($tr.Text -ne $newDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Walk every Design (theme) -> its SlideMaster -> the master's own Date
# placeholder, then every CustomLayout belonging to that master -> its
# Date placeholder (if the layout defines one).
for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $design = $p.Designs.Item($d)
    $master = $design.SlideMaster

    for ($i = 1; $i -le $master.Shapes.Count; $i++) {
        $shape = $master.Shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            Update-DateShape $shape
        }
    }

    for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
        $layout = $master.CustomLayouts.Item($L)
        for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
            $shape = $layout.Shapes.Item($i)
            if ($shape.Name -like "Date Placeholder*") {
                Update-DateShape $shape
            }
        }
    }
}

# Notes master has its own Date placeholder too.
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shape = $notesMaster.Shapes.Item($i)
    if ($shape.Name -like "Date Placeholder*") {
        Update-DateShape $shape
    }
}
